$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.213.18"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "1.861.25"
$ws.Range("E3").Value = "  -1.42%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'234.77"
$ws.Range("E5").Value = "  -1.41%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("E7").Value = "  -0.66%  "

$ws.Range("D8").Value = "'0.2826"
$ws.Range("E8").Value = "  -0.76%  "

$ws.Range("D9").Value = "'0.06526"
$ws.Range("E9").Value = "  -1.30%  "

$ws.Range("D10").Value = "'21.35"
$ws.Range("E10").Value = "  +3.83%  "

$ws.Range("D11").Value = "'0.07848"
$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").Value = "'97.14"
$ws.Range("E12").Value = "  -0.98%  "

$ws.Range("D13").Value = "1.866.46"
$ws.Range("E13").Value = "  -1.23%  "

$ws.Range("D14").Value = "'5.092"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").Value = "'0.6718"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("D16").Value = "'278.90"
$ws.Range("E16").Value = "  -2.39%  "

$ws.Range("D17").Value = "30.210.58"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").Value = "'5.510"
$ws.Range("E19").Value = "  +1.83%  "

$ws.Range("D20").Value = "'12.63"
$ws.Range("E20").Value = "  -0.36%  "

$ws.Range("D21").Value = "2.113.46"
$ws.Range("E21").Value = "  -1.14%  "

$ws.Range("D22").Value = "'0.000007270"
$ws.Range("E22").Value = "  -0.64%  "

$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "'6.139"
$ws.Range("E24").Value = "  -0.87%  "

$ws.Range("D25").Value = "'9.193"
$ws.Range("E25").Value = "  -2.43%  "

$ws.Range("D26").Value = "'164.75"
$ws.Range("E26").Value = "  -1.59%  "

$ws.Range("D27").Value = "'19.10"
$ws.Range("E27").Value = "  -1.03%  "

$ws.Range("D28").Value = "'1.919"
$ws.Range("E28").Value = "  -4.21%  "

$ws.Range("E29").Value = "  -0.60%  "

$ws.Range("D30").Value = "'0.09689"
$ws.Range("E30").Value = "  -0.87%  "

$ws.Range("D31").Value = "'4.419"
$ws.Range("E31").Value = "  +0.36%  "

$ws.Range("D32").Value = "'1.472"
$ws.Range("E32").Value = "  -1.17%  "

$ws.Range("D33").Value = "'4.072"
$ws.Range("E33").Value = "  -2.34%  "

$ws.Range("D34").Value = "'0.04686"
$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("D35").Value = "'1.113"
$ws.Range("E35").Value = "  +1.17%  "

$ws.Range("D36").Value = "'0.7036"
$ws.Range("E36").Value = "  -1.25%  "

$ws.Range("D37").Value = "'2.730"
$ws.Range("E37").Value = "  +0.72%  "

$ws.Range("D38").Value = "'0.01849"
$ws.Range("E38").Value = "  -1.80%  "

$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("D40").Value = "'6.238"
$ws.Range("E40").Value = "  -7.16%  "

$ws.Range("D41").Value = "'73.30"
$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("D42").Value = "'1.943"
$ws.Range("E42").Value = "  -1.99%  "

$ws.Range("D43").Value = "'0.8447"
$ws.Range("E43").Value = "  -3.20%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'104.03"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").Value = "'0.4154"
$ws.Range("E46").Value = "  -1.26%  "

$ws.Range("D47").Value = "'7.179"
$ws.Range("E47").Value = "  -1.80%  "

$ws.Range("D48").Value = "'936.29"
$ws.Range("E48").Value = "  -6.68%  "

$ws.Range("D49").Value = "'9.146"
$ws.Range("E49").Value = "  -0.65%  "

$ws.Range("D50").Value = "'33.99"
$ws.Range("E50").Value = "  -0.27%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1124"
$ws.Range("E51").Value = "  -2.84%  "
